$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.257.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.560.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.55%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.558.77"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.52%  "

$ws.Range("E9").Value = "  +3.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.12%  "

$ws.Range("E12").Value = "  +1.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.165.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.39%  "

$ws.Range("E14").Value = "  +4.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.558.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.351.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("E19").Value = "  +10.78%  "

$ws.Range("E20").Value = "  +1.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.609"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.706.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.63%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.34%  "

$ws.Range("E28").Value = "  +4.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("E31").Value = "  -0.13%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.160"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.558.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.33%  "

$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("E38").Value = "  +3.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.11%  "

$ws.Range("E40").Value = "  -0.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.64%  "

$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.898"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.41%  "

$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("E48").Value = "  +4.39%  "

$ws.Range("E49").Value = "  +5.57%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.32%  "
